$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for new columns I (I0) and J (IF).
# Copy the formatting (bold, border, center/top alignment) from the existing
# header cell H1 so the new header cells reuse the same cell style, then
# overwrite the copied value with the correct header text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Column I values (I0)
$ws.Range("I2").Value = 7
$ws.Range("I3").Value = 6
$ws.Range("I4").Value = 9
$ws.Range("I5").Value = 5
$ws.Range("I6").Value = 5
$ws.Range("I7").Value = 7

# Column J values (IF)
$ws.Range("J2").Value = 8
$ws.Range("J3").Value = 6
$ws.Range("J4").Value = 9
$ws.Range("J5").Value = 5
$ws.Range("J6").Value = 5
$ws.Range("J7").Value = 7
